# Auto-generated edit script: updates Leve profit-calculation values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (refreshed market-board pricing).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 7141.4736
$ws.Range("I28").Value = 1306.8889
$ws.Range("K28").Value = 1306.8889
$ws.Range("M28").Value = -821.8888999999999
$ws.Range("H112").Value = 2241.2246
$ws.Range("I112").Value = 969.6
$ws.Range("J112").Value = 2385.7273
$ws.Range("K112").Value = 2908.8
$ws.Range("L112").Value = 7157.1819
$ws.Range("M112").Value = -1800.8
$ws.Range("N112").Value = -9373.1819
$ws.Range("H138").Value = 3171.83
$ws.Range("I138").Value = 1575.3214
$ws.Range("K138").Value = 4725.9642
$ws.Range("M138").Value = 414.0357999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4777.1816
$ws.Range("I2").Value = 2021.75
$ws.Range("J2").Value = 12125
$ws.Range("K2").Value = 2021.75
$ws.Range("L2").Value = 12125
$ws.Range("M2").Value = -1908.75
$ws.Range("N2").Value = -12351
$ws.Range("H62").Value = 46689
$ws.Range("J62").Value = 46689
$ws.Range("L62").Value = 46689
$ws.Range("N62").Value = -47937
$ws.Range("H65").Value = 46689
$ws.Range("J65").Value = 46689
$ws.Range("L65").Value = 140067
$ws.Range("N65").Value = -146307
$ws.Range("H97").Value = 865.5
$ws.Range("I97").Value = 519.4
$ws.Range("J97").Value = 1607.1428
$ws.Range("K97").Value = 519.4
$ws.Range("L97").Value = 1607.1428
$ws.Range("M97").Value = -23.39999999999998
$ws.Range("N97").Value = -2599.1428
$ws.Range("H116").Value = 4777.1816
$ws.Range("I116").Value = 2021.75
$ws.Range("J116").Value = 12125
$ws.Range("K116").Value = 2021.75
$ws.Range("L116").Value = 12125
$ws.Range("M116").Value = 272.25
$ws.Range("N116").Value = -16713
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4777.1816
$ws.Range("I3").Value = 2021.75
$ws.Range("J3").Value = 12125
$ws.Range("K3").Value = 2021.75
$ws.Range("L3").Value = 12125
$ws.Range("M3").Value = -1907.75
$ws.Range("N3").Value = -12353

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20853014
$ws.Range("I31").Value = 55588856
$ws.Range("J31").Value = 11507.6
$ws.Range("K31").Value = 55588856
$ws.Range("L31").Value = 11507.6
$ws.Range("M31").Value = -55588561
$ws.Range("N31").Value = -12097.6
$ws.Range("H34").Value = 20853014
$ws.Range("I34").Value = 55588856
$ws.Range("J34").Value = 11507.6
$ws.Range("K34").Value = 55588856
$ws.Range("L34").Value = 11507.6
$ws.Range("M34").Value = -55588654
$ws.Range("N34").Value = -11911.6
$ws.Range("H141").Value = 713333.3
$ws.Range("J141").Value = 713333.3
$ws.Range("L141").Value = 713333.3
$ws.Range("N141").Value = -723693.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 104944.45
$ws.Range("J37").Value = 104944.45
$ws.Range("L37").Value = 314833.35
$ws.Range("N37").Value = -315057.35
$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 3000
$ws.Range("K88").Value = 9000
$ws.Range("M88").Value = -8572
$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 3000
$ws.Range("K91").Value = 9000
$ws.Range("M91").Value = -7518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 4209
$ws.Range("J6").Value = 4209
$ws.Range("L6").Value = 4209
$ws.Range("N6").Value = -4435
$ws.Range("H16").Value = 4209
$ws.Range("J16").Value = 4209
$ws.Range("L16").Value = 4209
$ws.Range("N16").Value = -4709
$ws.Range("H70").Value = 7503.5
$ws.Range("I70").Value = 7004.6
$ws.Range("K70").Value = 7004.6
$ws.Range("M70").Value = -6734.6
$ws.Range("H73").Value = 7503.5
$ws.Range("I73").Value = 7004.6
$ws.Range("K73").Value = 7004.6
$ws.Range("M73").Value = -6068.6
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H107").Value = 662.4706
$ws.Range("I107").Value = 172.4
$ws.Range("J107").Value = 866.6667
$ws.Range("K107").Value = 172.4
$ws.Range("L107").Value = 866.6667
$ws.Range("M107").Value = 1747.6
$ws.Range("N107").Value = -4706.6667
$ws.Range("H132").Value = 4870.5537
$ws.Range("I132").Value = 4778.7173
$ws.Range("J132").Value = 5293
$ws.Range("K132").Value = 14336.1519
$ws.Range("L132").Value = 15879
$ws.Range("M132").Value = -11806.1519
$ws.Range("N132").Value = -20939

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11417.368
$ws.Range("I7").Value = 8457.308000000001
$ws.Range("J7").Value = 17830.834
$ws.Range("K7").Value = 8457.308000000001
$ws.Range("L7").Value = 17830.834
$ws.Range("M7").Value = -8345.308000000001
$ws.Range("N7").Value = -18054.834
$ws.Range("H9").Value = 321
$ws.Range("I9").Value = 181.5
$ws.Range("K9").Value = 181.5
$ws.Range("M9").Value = 42.5
$ws.Range("H17").Value = 354503870
$ws.Range("I17").Value = 15035000
$ws.Range("K17").Value = 15035000
$ws.Range("M17").Value = -15034830
$ws.Range("H40").Value = 5229
$ws.Range("I40").Value = 4139.857
$ws.Range("K40").Value = 4139.857
$ws.Range("M40").Value = -4003.857
$ws.Range("H55").Value = 3128.16
$ws.Range("J55").Value = 8332.833000000001
$ws.Range("L55").Value = 8332.833000000001
$ws.Range("N55").Value = -8678.833000000001
$ws.Range("H68").Value = 4002.7778
$ws.Range("J68").Value = 5658.5
$ws.Range("L68").Value = 5658.5
$ws.Range("N68").Value = -7156.5
$ws.Range("H71").Value = 4002.7778
$ws.Range("J71").Value = 5658.5
$ws.Range("L71").Value = 28292.5
$ws.Range("N71").Value = -35780.5
$ws.Range("H82").Value = 2313.3333
$ws.Range("I82").Value = 1316.7
$ws.Range("K82").Value = 1316.7
$ws.Range("M82").Value = -955.7
$ws.Range("H85").Value = 2313.3333
$ws.Range("I85").Value = 1316.7
$ws.Range("K85").Value = 1316.7
$ws.Range("M85").Value = -68.70000000000005
$ws.Range("H122").Value = 5244.316
$ws.Range("I122").Value = 5068.923
$ws.Range("J122").Value = 5624.3335
$ws.Range("K122").Value = 15206.769
$ws.Range("L122").Value = 16873.0005
$ws.Range("M122").Value = -12756.769
$ws.Range("N122").Value = -21773.0005
$ws.Range("H126").Value = 11417.368
$ws.Range("I126").Value = 8457.308000000001
$ws.Range("J126").Value = 17830.834
$ws.Range("K126").Value = 25371.924
$ws.Range("L126").Value = 53492.50199999999
$ws.Range("M126").Value = -22901.924
$ws.Range("N126").Value = -58432.50199999999
$ws.Range("H132").Value = 2484.6287
$ws.Range("I132").Value = 2264
$ws.Range("J132").Value = 3551
$ws.Range("K132").Value = 6792
$ws.Range("L132").Value = 10653
$ws.Range("M132").Value = -4262
$ws.Range("N132").Value = -15713
$ws.Range("H136").Value = 35723904
$ws.Range("I136").Value = 83342290
$ws.Range("K136").Value = 250026870
$ws.Range("M136").Value = -250024320

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 19999
$ws.Range("J32").Value = 19999
$ws.Range("L32").Value = 19999
$ws.Range("N32").Value = -20633
$ws.Range("H68").Value = 9999
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 9999
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
